$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C3: presente -> ausente
$ws.Range("C3").Value = "ausente"

# Add new row 10
$ws.Range("A10").Value = 1234
$ws.Range("B10").Value = "fdsafdasfs"
$ws.Range("C10").Value = "ausente"
$ws.Range("D10").Value = 9
$ws.Range("E10").Value = "B"
